# Add an "Electrode Locations" column to the main dataframe and sort the
# data rows by electrode location (letter, then numeric position), A1-O15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (data starts at row 2; row 1 is header).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

# Read existing File Name / Unnormalized P_max pairs into memory.
$records = @()
for ($r = 2; $r -le $lastRow; $r++) {
    $fileName = $ws.Cells.Item($r, 1).Value2
    if ($null -eq $fileName -or $fileName -eq "") { continue }
    $value = $ws.Cells.Item($r, 2).Value2

    # Electrode location = filename text before the first underscore,
    # e.g. "A11_bipolar_10V_100kHz.txt" -> "A11"
    $location = $fileName.Split("_")[0]

    # Split location into its leading letters and trailing number so the
    # rows can be sorted "A1, A2, ... A15, C1, C3, ... O11" instead of
    # plain alphabetical text order.
    if ($location -match '^([A-Za-z]+)(\d+)$') {
        $letters = $matches[1]
        $number = [int]$matches[2]
    } else {
        $letters = $location
        $number = 0
    }

    $records += @{
        FileName = $fileName
        Value = $value
        Location = $location
        Letters = $letters
        Number = $number
    }
}

# Sort by letter prefix, then by numeric suffix (A2, A3, A5 ... A15, C1, C3 ...).
# Sort-Object is stable, so sorting by the secondary key first and then by
# the primary key reproduces a combined (Letters, Number) ordering.
$byNumber = $records | Sort-Object -Property Number
$sorted = $byNumber | Sort-Object -Property Letters

# Header for the new column - reuse the same formatting as the other headers.
$ws.Cells.Item(1, 1).Copy()
$ws.Cells.Item(1, 3).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(1, 3).Value2 = "Electrode Locations"

# Write the sorted rows back out, filling columns A, B and the new C.
$r = 2
foreach ($rec in $sorted) {
    $ws.Cells.Item($r, 1).Value2 = $rec.FileName
    $ws.Cells.Item($r, 2).Value2 = $rec.Value
    $ws.Cells.Item($r, 3).Value2 = $rec.Location
    $r++
}

